# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload"
#
# Fills in the two placeholder rows (22 & 23) of the JAN-22 daily-tracker
# sheet with the 12th day's entries, and moves the active-cell selection
# to D23 (where the author was last typing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-22")

# --- Row 23 first, then row 22 --------------------------------------
# (sets the new shared-string table entries in the same order Excel
#  wrote them: "Also, we are working..." before "Enhancing the
#  knowledge...")
$ws.Range("D23").Value = "2. Also, we are working date filtering with new logic to obtain previous month data if vendor ledger is having three months data and `nit is work in progress for the GL"
$ws.Range("D23").WrapText = $true
$ws.Range("E23").Value = 0.3
$ws.Range("E23").NumberFormat = "0%"
$ws.Range("F23").Value = "WIP"

$ws.Range("A22").Value = 12
$ws.Range("B22").Value = 44582
$ws.Range("B22").NumberFormat = "m/d/yy"
$ws.Range("C22").Value = "RPA RLOGIC"
$ws.Range("D22").Value = "1. Enhancing the knowledge of Openpyxl and Pandas frameworks towards the RPA automation to handle excel manipulation and involved in the `ncalculation activities "
$ws.Range("D22").WrapText = $true
$ws.Range("E22").Value = 1
$ws.Range("E22").NumberFormat = "0%"
$ws.Range("F22").Value = "Completed"

# Row heights grow to fit the wrapped two-line comments.
$ws.Rows.Item(22).RowHeight = 28.8
$ws.Rows.Item(23).RowHeight = 28.8

# Leave the selection where the author left off.
[void]$ws.Range("D23").Select()
